$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "1.000", "0.4700") that
# Excel would otherwise auto-convert to a Number on assignment, stripping the
# literal formatting. Force the whole Price column to Text first, write the
# values, then clear the temporary number format again so the saved cells carry
# no extra style (matching the original un-styled cells).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.293.66"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.869.59"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "235.05"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "0.06594"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "0.07960"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "96.80"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "1.879.97"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "0.6972"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "5.115"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "268.58"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "30.343.65"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "14.16"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").Value = "0.000007789"
$ws.Range("E19").Value = "  +5.93%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "2.121.25"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "5.271"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").Value = "6.218"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "9.384"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").Value = "167.48"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").Value = "0.09919"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "4.336"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "1.459"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "4.055"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "0.04721"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "0.7036"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "0.01875"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +6.91%  "
$ws.Range("D40").Value = "6.266"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "72.35"
$ws.Range("E41").Value = "  -4.48%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4177"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8422"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "102.92"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "7.117"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "9.190"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "919.16"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").Value = "34.61"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "0.05683"
$ws.Range("E51").Value = "  +0.58%  "

$priceRange.ClearFormats()
